$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'63.756.65"
$ws.Range("E2").Value = "'  -1.59%  "
$ws.Range("D3").Value = "'3.403.35"
$ws.Range("E3").Value = "'  -1.04%  "
$ws.Range("E4").Value = "'  +0.01%  "
$ws.Range("D5").Value = "'569.84"
$ws.Range("E5").Value = "'  -0.61%  "
$ws.Range("D6").Value = "'161.62"
$ws.Range("E6").Value = "'  +1.61%  "
$ws.Range("E7").Value = "'  +0.00%  "
$ws.Range("D8").Value = "'3.404.01"
$ws.Range("E8").Value = "'  -1.10%  "
$ws.Range("D9").Value = "'0.550"
$ws.Range("E9").Value = "'  -5.33%  "
$ws.Range("D10").Value = "'7.27"
$ws.Range("E10").Value = "'  +1.02%  "
$ws.Range("D11").Value = "'0.119"
$ws.Range("E11").Value = "'  -2.18%  "
$ws.Range("D12").Value = "'0.422"
$ws.Range("E12").Value = "'  -4.85%  "
$ws.Range("D13").Value = "'3.994.31"
$ws.Range("E13").Value = "'  -0.97%  "
$ws.Range("E14").Value = "'  +1.08%  "
$ws.Range("D15").Value = "'26.87"
$ws.Range("E15").Value = "'  -2.81%  "
$ws.Range("E16").Value = "'  -4.37%  "
$ws.Range("D17").Value = "'63.815.67"
$ws.Range("E17").Value = "'  -1.53%  "
$ws.Range("D18").Value = "'3.426.58"
$ws.Range("E18").Value = "'  -0.45%  "
$ws.Range("D19").Value = "'6.10"
$ws.Range("E19").Value = "'  -3.85%  "
$ws.Range("D20").Value = "'13.53"
$ws.Range("E20").Value = "'  -2.62%  "
$ws.Range("D21").Value = "'375.50"
$ws.Range("E21").Value = "'  -1.49%  "
$ws.Range("D22").Value = "'7.74"
$ws.Range("E22").Value = "'  -2.87%  "
$ws.Range("D23").Value = "'0.998"
$ws.Range("E23").Value = "'  -0.20%  "
$ws.Range("D24").Value = "'70.94"
$ws.Range("E24").Value = "'  -1.63%  "
$ws.Range("D25").Value = "'0.514"
$ws.Range("E25").Value = "'  -6.15%  "
$ws.Range("E26").Value = "'  -3.25%  "
$ws.Range("D27").Value = "'9.49"
$ws.Range("E27").Value = "'  -3.61%  "
$ws.Range("E28").Value = "'  +0.23%  "
$ws.Range("E29").Value = "'  -0.12%  "
$ws.Range("D30").Value = "'6.07"
$ws.Range("E30").Value = "'  -0.38%  "
$ws.Range("E31").Value = "'  -6.36%  "
$ws.Range("D32").Value = "'2.00"
$ws.Range("E32").Value = "'  -0.20%  "
$ws.Range("D33").Value = "'22.84"
$ws.Range("D34").Value = "'7.06"
$ws.Range("E34").Value = "'  +0.82%  "
$ws.Range("E35").Value = "'  -4.61%  "
$ws.Range("D36").Value = "'159.81"
$ws.Range("E36").Value = "'  -0.82%  "
$ws.Range("D37").Value = "'0.857"
$ws.Range("E37").Value = "'  +10.12%  "
$ws.Range("D38").Value = "'1.80"
$ws.Range("E38").Value = "'  -4.57%  "
$ws.Range("D39").Value = "'0.0724"
$ws.Range("E39").Value = "'  -2.77%  "
$ws.Range("D40").Value = "'2.771.01"
$ws.Range("E40").Value = "'  -4.14%  "
$ws.Range("D41").Value = "'25.65"
$ws.Range("E41").Value = "'  -2.18%  "
$ws.Range("D42").Value = "'42.70"
$ws.Range("E42").Value = "'  -0.51%  "
$ws.Range("D43").Value = "'6.41"
$ws.Range("E43").Value = "'  -3.40%  "
$ws.Range("D44").Value = "'25.99"
$ws.Range("E44").Value = "'  +0.12%  "
$ws.Range("D45").Value = "'4.38"
$ws.Range("E45").Value = "'  -3.31%  "
$ws.Range("D46").Value = "'0.0305"
$ws.Range("E46").Value = "'  -3.49%  "
$ws.Range("D47").Value = "'2.39"
$ws.Range("E47").Value = "'  +5.42%  "
$ws.Range("D48").Value = "'328.84"
$ws.Range("E48").Value = "'  +4.10%  "
$ws.Range("E49").Value = "'  -4.10%  "
$ws.Range("D50").Value = "'6.28"
$ws.Range("E50").Value = "'  -3.48%  "
$ws.Range("E51").Value = "'  -3.08%  "
